# Update the Metadata sheet of the StructureDefinition workbook.
# Changes:
#  1. Title value becomes "Extension String Fundamento Priorizacion"
#     (previously reused the Name value "ExtensionStringFundamentoPriorizacion").
#  2. Date value is bumped to the new publication timestamp.
#  3. Context value is simplified from a full URL-based element context
#     to the short "element:ServiceRequest" form.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Cells.Item(5, 2).Value = "Extension String Fundamento Priorizacion"
$ws.Cells.Item(8, 2).Value = "2024-07-15T11:25:06-04:00"
$ws.Cells.Item(21, 2).Value = "element:ServiceRequest"
